$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Drop the "Meta description: ..." paragraph that used to follow
#    the title (Heading1) at the top of the document.
# ------------------------------------------------------------------
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# ------------------------------------------------------------------
# 2. At the very end of the document, the italic "Create a feature
#    image..." paragraph is replaced by two paragraphs:
#      - a new bold paragraph repeating the page title
#      - the former "Meta description" text (minus the "Meta
#        description" label) now in italics
# ------------------------------------------------------------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

$newXml = @"
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Fruits &amp; 777's for Free - Review of Classic Slot Game</w:t></w:r></w:p><w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Read our review of Fruits &amp; 777's and play for free. Learn about the RTP range, betting options, and what we like and don't like about this classic slot game.</w:t></w:r></w:p>
"@

$lastPara.Range.InsertXML($newXml)
